$wb = $excel.ActiveWorkbook

# Overview sheet: update "Latest HO Xliff Generate Date" for first row
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G2").Value = "2017-02-09 16:12:29"

# zh-cn sheet: update Correspond Handoff Datetime and Correspond Handback DateTime
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2017-02-09 16:12:10"
$wsZhCn.Range("L2").Value = "2017-02-09 16:12:58"

# de-de sheet: update Correspond Handoff Datetime and Correspond Handback DateTime
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2017-02-09 16:12:29"
$wsDeDe.Range("L2").Value = "2017-02-09 16:13:22"
